# C5-PowerPoint.pptx edit
#
# 1) Change the table on slide 6 from the custom "Table_0" style
#    ({AABE5E44-460D-4C2D-A33F-57434E15EDED}) to the built-in PowerPoint
#    table style {F9B3B137-424E-4887-AAD9-D165240CE887}.
#
# 2) Swap the two theme colour schemes: the deck's main theme (applied to
#    the slide master, currently the green "Integral" palette) is changed
#    back to the default "Office Theme" colour palette.

$p = $ppt.ActivePresentation

# --- 1. Table style -------------------------------------------------------

$targetStyleId = "{F9B3B137-424E-4887-AAD9-D165240CE887}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($targetStyleId)
        }
    }
}

# --- 2. Theme colours ------------------------------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink (standard theme colour
# order exposed via ThemeColorScheme.Colors(1..12)). These are the
# "Office Theme" values that the deck's theme is being restored to.

$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le 12; $i++) {
    $themeColors.Colors($i).RGB = $officeColors[$i - 1]
}
